$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 162, shifting existing rows 162-191 down to 163-192
$ws.Rows.Item(162).Insert()

# Populate the new row 162 with the new weekly record
$ws.Cells.Item(162, 1).Value = 4
$ws.Cells.Item(162, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(162, 3).Value = "Los Lagos"
$ws.Cells.Item(162, 4).Value = (Get-Date -Year 2021 -Month 10 -Day 5 -Hour 0 -Minute 0 -Second 0)
$ws.Cells.Item(162, 5).Value = 10
$ws.Cells.Item(162, 6).Value = 100112008
$ws.Cells.Item(162, 7).Value = "Coliflor"
$ws.Cells.Item(162, 8).Value = "Sin especificar"
$ws.Cells.Item(162, 9).Value = "Primera"
$ws.Cells.Item(162, 10).Value = 1400
$ws.Cells.Item(162, 11).Value = 1000
$ws.Cells.Item(162, 12).Value = 1100
$ws.Cells.Item(162, 13).Value = 1050
$ws.Cells.Item(162, 14).Value = "`$/unidad"
$ws.Cells.Item(162, 15).Value = "Región Metropolitana"
$ws.Cells.Item(162, 16).Value = 1050
$ws.Cells.Item(162, 17).Value = 1
$ws.Cells.Item(162, 18).Value = "Hortaliza"
